$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to text format so numeric-looking strings
# (e.g. "1.00", "28.03") are stored as text, matching the source data,
# then clear the formatting override so no new cell style lingers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.442.37"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.648.84"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "597.77"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").Value = "159.00"
$ws.Range("E6").Value = "  +2.68%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "28.03"
$ws.Range("D14").Value = "3.124.24"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "0.0000187"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("D16").Value = "68.305.22"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "2.628.37"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "11.42"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "364.15"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").Value = "4.41"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "4.78"
$ws.Range("E22").Value = "  -2.72%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "74.41"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "9.80"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").Value = "2.780.82"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "0.0000103"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "560.53"
$ws.Range("E30").Value = "  -2.19%  "
$ws.Range("D31").Value = "8.06"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "1.65"
$ws.Range("E34").Value = "  +3.64%  "
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "160.66"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").Value = "1.87"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("E43").Value = "  -4.44%  "
$ws.Range("D45").Value = "158.22"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").Value = "22.06"
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("D48").Value = "1.69"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").Value = "0.0773"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("E51").Value = "  -0.93%  "

$ws.Range("D2:D51").ClearFormats()
